{"js": "// Update \"21 years\" -> \"15+ years\" in the professional summary,\n// and remove the EDUCATION section (heading + its two entries).\n\n// 1. Update experience years in the professional summary paragraph.\nconst summaryResults = context.document.body.search(\n  \"21 years of experience in full-stack development, data engineering, and scalable web applications.\",\n  { matchCase: true }\n);\nsummaryResults.load(\"text\");\nawait context.sync();\n\nif (summaryResults.items.length > 0) {\n  summaryResults.items[0].insertText(\n    \"15+ years of experience in full-stack development, data engineering, and scalable web applications.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 2. Remove the EDUCATION section: the \"EDUCATION\" heading paragraph and\n// the two degree entries that follow it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst targets = [\n  \"EDUCATION\",\n  \"Master of Arts in Political Science - University of California, Berkeley\",\n  \"Bachelor of Arts in Political Science - University of California, Berkeley\"\n];\n\nconst paragraphsToDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text.trim();\n  if (targets.indexOf(text) !== -1) {\n    paragraphsToDelete.push(paragraph);\n  }\n}\n\nfor (const paragraph of paragraphsToDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Update \"21 years\" -> \"15+ years\" in the professional summary,\n# and remove the EDUCATION section (heading + its two entries).\n\n$d = $word.ActiveDocument\n\n# 1. Update experience years in the professional summary paragraph.\n$find = $d.Content.Find\n$find.Text = \"21 years of experience in full-stack development, data engineering, and scalable web applications.\"\n$find.Replacement.Text = \"15+ years of experience in full-stack development, data engineering, and scalable web applications.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Remove the EDUCATION section: the \"EDUCATION\" heading paragraph and\n# the two degree entries that follow it.\n$targets = @(\n    \"EDUCATION\",\n    \"Master of Arts in Political Science - University of California, Berkeley\",\n    \"Bachelor of Arts in Political Science - University of California, Berkeley\"\n)\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($targets -contains $t) {\n        $p.Range.Delete()\n    }\n}\n"}
